$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '56.914.89'
Set-TextValue 'E2' '  -0.58%  '
Set-TextValue 'D3' '2.317.93'
Set-TextValue 'E3' '  -1.61%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.26%  '
Set-TextValue 'D5' '530.03'
Set-TextValue 'E5' '  +2.16%  '
Set-TextValue 'E6' '  -2.15%  '
Set-TextValue 'D7' '0.996'
Set-TextValue 'E7' '  -0.15%  '
Set-TextValue 'E8' '  -0.87%  '
Set-TextValue 'D9' '2.341.05'
Set-TextValue 'E9' '  -1.36%  '
Set-TextValue 'D10' '0.101'
Set-TextValue 'E10' '  -1.28%  '
Set-TextValue 'E11' '  +0.02%  '
Set-TextValue 'E12' '  -3.12%  '
Set-TextValue 'E13' '  +1.51%  '
Set-TextValue 'D14' '2.735.76'
Set-TextValue 'E14' '  -1.53%  '
Set-TextValue 'D15' '23.41'
Set-TextValue 'E15' '  -3.70%  '
Set-TextValue 'D16' '56.965.15'
Set-TextValue 'E16' '  -0.50%  '
Set-TextValue 'E17' '  -1.77%  '
Set-TextValue 'D18' '2.337.48'
Set-TextValue 'E18' '  -0.98%  '
Set-TextValue 'D19' '335.89'
Set-TextValue 'D20' '10.40'
Set-TextValue 'E20' '  -1.57%  '
Set-TextValue 'E21' '  -1.83%  '
Set-TextValue 'D22' '6.83'
Set-TextValue 'E22' '  +1.89%  '
Set-TextValue 'D23' '0.998'
Set-TextValue 'E23' '  +0.06%  '
Set-TextValue 'D24' '61.68'
Set-TextValue 'E24' '  +0.58%  '
Set-TextValue 'E25' '  +1.17%  '
Set-TextValue 'D26' '8.68'
Set-TextValue 'E26' '  -2.42%  '
Set-TextValue 'E27' '  -0.27%  '
Set-TextValue 'D28' '1.36'
Set-TextValue 'E28' '  +3.05%  '
Set-TextValue 'D29' '172.33'
Set-TextValue 'E29' '  +3.02%  '
Set-TextValue 'D31' '0.0₃0724'
Set-TextValue 'E31' '  -2.38%  '
Set-TextValue 'D32' '6.10'
Set-TextValue 'E32' '  -2.87%  '
Set-TextValue 'D33' '18.46'
Set-TextValue 'E33' '  -0.54%  '
Set-TextValue 'D34' '0.999'
Set-TextValue 'E34' '  -0.05%  '
Set-TextValue 'D35' '0.992'
Set-TextValue 'E35' '  -0.29%  '
Set-TextValue 'E36' '  -3.32%  '
Set-TextValue 'D37' '0.929'
Set-TextValue 'E37' '  +0.77%  '
Set-TextValue 'D38' '3.98'
Set-TextValue 'E38' '  -0.84%  '
Set-TextValue 'D39' '39.19'
Set-TextValue 'E39' '  +0.89%  '
Set-TextValue 'E40' '  -1.91%  '
Set-TextValue 'B41' 'Aave'
Set-TextValue 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '148.66'
Set-TextValue 'E41' '  -0.44%  '
Set-TextValue 'B42' 'RenderToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D42' '5.59'
Set-TextValue 'E42' '  +6.69%  '
Set-TextValue 'E43' '  -2.93%  '
Set-TextValue 'D44' '3.60'
Set-TextValue 'E44' '  -1.29%  '
Set-TextValue 'D45' '281.41'
Set-TextValue 'E45' '  -1.70%  '
Set-TextValue 'D46' '0.0932'
Set-TextValue 'E46' '  -0.86%  '
Set-TextValue 'D47' '0.0500'
Set-TextValue 'E47' '  -1.66%  '
Set-TextValue 'D48' '18.81'
Set-TextValue 'E48' '  +3.40%  '
Set-TextValue 'D49' '0.558'
Set-TextValue 'E49' '  -1.21%  '
Set-TextValue 'D50' '0.0216'
Set-TextValue 'E50' '  -1.24%  '
Set-TextValue 'B51' 'Polygon'
Set-TextValue 'C51' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D51' '0.382'
Set-TextValue 'E51' '  +5.55%  '
